$d = $word.ActiveDocument
$last = $d.Paragraphs.Last
$secondLast = $last.Previous()
$rng = $secondLast.Range
$rng.Collapse(1)
$frag = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:spacing w:before="35"/>
        <w:ind w:left="1420"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:spacing w:before="35"/>
        <w:ind w:left="1420"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:spacing w:before="35"/>
        <w:ind w:left="1420"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:spacing w:before="35"/>
      </w:pPr>
      <w:r>
        <w:t>Tabel</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:spacing w:val="-2"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Laporan</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:spacing w:val="-3"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Hasil</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:spacing w:val="-1"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Sprint</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:spacing w:val="-3"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Review</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
        <w:spacing w:before="9"/>
        <w:rPr>
          <w:sz w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="-5" w:type="dxa"/>
        <w:tblBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:insideH w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:insideV w:val="single" w:sz="4" w:space="0" w:color="000000"/>
        </w:tblBorders>
        <w:tblLayout w:type="fixed"/>
        <w:tblCellMar>
          <w:left w:w="0" w:type="dxa"/>
          <w:right w:w="0" w:type="dxa"/>
        </w:tblCellMar>
        <w:tblLook w:val="01E0" w:firstRow="1" w:lastRow="1" w:firstColumn="1" w:lastColumn="1" w:noHBand="0" w:noVBand="0"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="3115"/>
        <w:gridCol w:w="3117"/>
        <w:gridCol w:w="3117"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="556"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9349" w:type="dxa"/>
            <w:gridSpan w:val="3"/>
            <w:tcBorders>
              <w:top w:val="nil"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="7F7F7F" w:themeFill="text1" w:themeFillTint="80"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="9"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FFFFFF"/>
              </w:rPr>
              <w:t>Sprint</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:color w:val="FFFFFF"/>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:color w:val="FFFFFF"/>
              </w:rPr>
              <w:t>Review</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="9" w:line="249" w:lineRule="exact"/>
            </w:pPr>
            <w:r>
              <w:t>Sprint</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>1</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="268"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3115" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="248" w:lineRule="exact"/>
            </w:pPr>
            <w:r>
              <w:t>Completed</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Task</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="248" w:lineRule="exact"/>
              <w:ind w:left="108"/>
            </w:pPr>
            <w:r>
              <w:t>Incompleted</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-3"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Task</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="248" w:lineRule="exact"/>
              <w:ind w:left="108"/>
            </w:pPr>
            <w:r>
              <w:t>What’s</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-1"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Next</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="1074"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3115" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:right="550"/>
            </w:pPr>
            <w:r>
              <w:t>Task 1: Membuat Tampilan</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-47"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Figma</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="270" w:lineRule="atLeast"/>
              <w:ind w:right="498"/>
            </w:pPr>
            <w:r>
              <w:t>Task</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-7"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>2:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-4"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Membuat</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-6"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Prototype</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-46"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Figma</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="8"/>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:t>-</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="108"/>
            </w:pPr>
            <w:r>
              <w:t>Sprint</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-3"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>2:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Sprintgoal1</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="266"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3115" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="246" w:lineRule="exact"/>
            </w:pPr>
            <w:r>
              <w:t>Sprint</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>2</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman"/>
                <w:sz w:val="18"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman"/>
                <w:sz w:val="18"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="1074"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3115" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:right="550"/>
            </w:pPr>
            <w:r>
              <w:t>Task 1: Membuat Tampilan</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-47"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Figma</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="267" w:lineRule="exact"/>
            </w:pPr>
            <w:r>
              <w:t>Task</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-4"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>2: Membuat</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:spacing w:val="-2"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>Prototype</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:line="252" w:lineRule="exact"/>
            </w:pPr>
            <w:r>
              <w:t>Figma</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="8"/>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:t>-</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Sprint 3: Sprintgoal2</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="1074"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3115" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:ind w:right="550"/>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Sprint3</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="21"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="21"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="21"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>-</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3117" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="TableParagraph"/>
              <w:spacing w:before="11"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:sz w:val="21"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="21"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Completed</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>

'@
$rng.InsertXML($frag)
Write-Host "Inserted report table."
